# Commit: "finished handling panel and defaults, begin working on coupons"
#
# Changes applied:
#  1. The coupon tables' "blade width (leave blank)" header is renamed to
#     "blade width (leave blank if same)" (cells E10 and D19 on the
#     "Example" sheet).
#  2. Columns D and E are widened (to fit the longer new header text) to a
#     matching width.
#  3. The active selection on the "Example" sheet moves from C3:D3 down to
#     A21 (the row right below the last used row), reflecting that the
#     author finished the "Panel"/"Defaults" sections and is starting on
#     the coupons sections below.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")
$ws.Activate()

# 1. Update the "blade width" column headers for both coupon tables.
$ws.Range("E10").Value = "blade width (leave blank if same)"
$ws.Range("D19").Value = "blade width (leave blank if same)"

# 2. Widen columns D:E to accommodate the longer text.
$ws.Columns.Item(4).ColumnWidth = 25
$ws.Columns.Item(5).ColumnWidth = 25

# 3. Move the selection to A21, ready for the next block of work.
$ws.Range("A21").Select() | Out-Null
